$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 66.47695399999999
$ws.Range("H2").Value = 199.430862
$ws.Range("I2").Value = 0.04311983106164722
$ws.Range("J2").Value = 0.04311983106164721
$ws.Range("M2").Value = 0.08241233333333334
$ws.Range("Q2").Value = 5.478520892032666
$ws.Range("R2").Value = 49.306688028294
$ws.Range("S2").Value = 0.002038276834540459
$ws.Range("T2").Value = 0.002038276834540459
$ws.Range("G3").Value = 66.47695399999999
$ws.Range("H3").Value = 199.430862
$ws.Range("I3").Value = 0.04311983106164722
$ws.Range("J3").Value = 0.04311983106164721
$ws.Range("Q3").Value = 110.419816040896
$ws.Range("R3").Value = 993.7783443680639
$ws.Range("S3").Value = 0.04108155422710676
$ws.Range("T3").Value = 0.04108155422710675
$ws.Range("I4").Value = 0.8830494168872806
$ws.Range("J4").Value = 0.8830494168872804
$ws.Range("M4").Value = 0.08241233333333334
$ws.Range("S4").Value = 0.04174179550060247
$ws.Range("T4").Value = 0.04174179550060246
$ws.Range("I5").Value = 0.8830494168872806
$ws.Range("J5").Value = 0.8830494168872804
$ws.Range("S5").Value = 0.8413076213866781
$ws.Range("T5").Value = 0.8413076213866779
$ws.Range("G6").Value = 44.831112
$ws.Range("H6").Value = 134.493336
$ws.Range("I6").Value = 0.02907940059566787
$ws.Range("J6").Value = 0.02907940059566786
$ws.Range("M6").Value = 0.08241233333333334
$ws.Range("Q6").Value = 3.694636545848
$ws.Range("R6").Value = 33.251728912632
$ws.Range("S6").Value = 0.0013745848983437
$ws.Range("T6").Value = 0.001374584898343699
$ws.Range("G7").Value = 44.831112
$ws.Range("H7").Value = 134.493336
$ws.Range("I7").Value = 0.02907940059566787
$ws.Range("J7").Value = 0.02907940059566786
$ws.Range("Q7").Value = 74.465552978688
$ws.Range("R7").Value = 670.189976808192
$ws.Range("S7").Value = 0.02770481569732417
$ws.Range("T7").Value = 0.02770481569732416
$ws.Range("G8").Value = 52.83062100000001
$ws.Range("H8").Value = 158.491863
$ws.Range("I8").Value = 0.0342682285413064
$ws.Range("J8").Value = 0.03426822854130639
$ws.Range("M8").Value = 0.08241233333333334
$ws.Range("Q8").Value = 4.353894748059001
$ws.Range("R8").Value = 39.18505273253101
$ws.Range("S8").Value = 0.001619861086575759
$ws.Range("T8").Value = 0.001619861086575758
$ws.Range("G9").Value = 52.83062100000001
$ws.Range("H9").Value = 158.491863
$ws.Range("I9").Value = 0.0342682285413064
$ws.Range("J9").Value = 0.03426822854130639
$ws.Range("Q9").Value = 87.75292941590402
$ws.Range("R9").Value = 789.7763647431361
$ws.Range("S9").Value = 0.03264836745473065
$ws.Range("T9").Value = 0.03264836745473063
$ws.Range("G10").Value = 16.16161433333333
$ws.Range("H10").Value = 48.484843
$ws.Range("I10").Value = 0.01048312291409786
$ws.Range("J10").Value = 0.01048312291409786
$ws.Range("M10").Value = 0.08241233333333334
$ws.Range("Q10").Value = 1.331916347643444
$ws.Range("R10").Value = 11.987247128791
$ws.Range("S10").Value = 0.0004955378085525757
$ws.Range("T10").Value = 0.0004955378085525756
$ws.Range("G11").Value = 16.16161433333333
$ws.Range("H11").Value = 48.484843
$ws.Range("I11").Value = 0.01048312291409786
$ws.Range("J11").Value = 0.01048312291409786
$ws.Range("Q11").Value = 26.84482928641066
$ws.Range("R11").Value = 241.603463577696
$ws.Range("S11").Value = 0.00998758510554529
$ws.Range("T11").Value = 0.009987585105545286
